$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")
$ws.Select()
